$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 484.07693
$ws.Range("J17").Value = 484.07693
$ws.Range("L17").Value = 1452.23079
$ws.Range("N17").Value = -1788.23079

$ws.Range("H53").Value = 455.14285
$ws.Range("I53").Value = 114.833336
$ws.Range("J53").Value = 710.375
$ws.Range("K53").Value = 114.833336
$ws.Range("L53").Value = 710.375
$ws.Range("M53").Value = 522.166664
$ws.Range("N53").Value = -1984.375

$ws.Range("H81").Value = 78629.25
$ws.Range("J81").Value = 78629.25
$ws.Range("L81").Value = 78629.25
$ws.Range("N81").Value = -80625.25

$ws.Range("H84").Value = 78629.25
$ws.Range("J84").Value = 78629.25
$ws.Range("L84").Value = 235887.75
$ws.Range("N84").Value = -245871.75

$ws.Range("H96").Value = 1679.625
$ws.Range("I96").Value = 1662
$ws.Range("J96").Value = 1697.25
$ws.Range("K96").Value = 4986
$ws.Range("L96").Value = 5091.75
$ws.Range("M96").Value = -3613
$ws.Range("N96").Value = -7837.75

$ws.Range("H100").Value = 3257.6667
$ws.Range("I100").Value = 1497.2858
$ws.Range("J100").Value = 5722.2
$ws.Range("K100").Value = 1497.2858
$ws.Range("L100").Value = 5722.2
$ws.Range("M100").Value = -956.2858000000001
$ws.Range("N100").Value = -6804.2

$ws.Range("H106").Value = 4999.2
$ws.Range("I106").Value = 3499.5
$ws.Range("K106").Value = 3499.5
$ws.Range("M106").Value = -2868.5

$ws.Range("H111").Value = 3730.8
$ws.Range("I111").Value = 2341.3333
$ws.Range("J111").Value = 5815
$ws.Range("K111").Value = 7023.999899999999
$ws.Range("L111").Value = 17445
$ws.Range("M111").Value = -3956.999899999999
$ws.Range("N111").Value = -23579

$ws.Range("H128").Value = 76516.664
$ws.Range("J128").Value = 76516.664
$ws.Range("L128").Value = 76516.664
$ws.Range("N128").Value = -86476.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10874293
$ws.Range("I32").Value = 11368573
$ws.Range("J32").Value = 149
$ws.Range("K32").Value = 11368573
$ws.Range("L32").Value = 149
$ws.Range("M32").Value = -11368286
$ws.Range("N32").Value = -723

$ws.Range("H45").Value = 2234.7646
$ws.Range("I45").Value = 1750.1666
$ws.Range("K45").Value = 1750.1666
$ws.Range("M45").Value = -1373.1666

$ws.Range("H110").Value = 8708.777
$ws.Range("I110").Value = 6197
$ws.Range("K110").Value = 6197
$ws.Range("M110").Value = -4152

$ws.Range("H122").Value = 1882.6666
$ws.Range("I122").Value = 1598
$ws.Range("K122").Value = 4794
$ws.Range("M122").Value = -2344

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 121000
$ws.Range("J62").Value = 122000
$ws.Range("L62").Value = 122000
$ws.Range("N62").Value = -123372

$ws.Range("H65").Value = 121000
$ws.Range("J65").Value = 122000
$ws.Range("L65").Value = 366000
$ws.Range("N65").Value = -372864

$ws.Range("H94").Value = 1493.8823
$ws.Range("I94").Value = 1426.4667
$ws.Range("K94").Value = 1426.4667
$ws.Range("M94").Value = -975.4666999999999

$ws.Range("H105").Value = 2062.5
$ws.Range("I105").Value = 1783.3334
$ws.Range("K105").Value = 1783.3334
$ws.Range("M105").Value = -36.33339999999998

$ws.Range("H107").Value = 1679.6428
$ws.Range("I107").Value = 1458.75
$ws.Range("K107").Value = 1458.75
$ws.Range("M107").Value = 461.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 487.5
$ws.Range("I16").Value = 550
$ws.Range("K16").Value = 550
$ws.Range("M16").Value = -263

$ws.Range("H31").Value = 904166.7
$ws.Range("I31").Value = 12049.786
$ws.Range("K31").Value = 12049.786
$ws.Range("M31").Value = -11754.786

$ws.Range("H34").Value = 904166.7
$ws.Range("I34").Value = 12049.786
$ws.Range("K34").Value = 12049.786
$ws.Range("M34").Value = -11847.786

$ws.Range("H105").Value = 2117.375
$ws.Range("I105").Value = 1732.25
$ws.Range("K105").Value = 1732.25
$ws.Range("M105").Value = 14.75

$ws.Range("H113").Value = 487.5
$ws.Range("I113").Value = 550
$ws.Range("K113").Value = 550
$ws.Range("M113").Value = 1620

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1892.7693
$ws.Range("I122").Value = 975.75
$ws.Range("J122").Value = 2300.3333
$ws.Range("K122").Value = 8781.75
$ws.Range("L122").Value = 20702.9997
$ws.Range("M122").Value = -6331.75
$ws.Range("N122").Value = -25602.9997

$ws.Range("H128").Value = 417999.66
$ws.Range("I128").Value = 417999.66
$ws.Range("K128").Value = 1253998.98
$ws.Range("M128").Value = -1249018.98

$ws.Range("H131").Value = 1959.5834
$ws.Range("I131").Value = 2043
$ws.Range("K131").Value = 6129
$ws.Range("M131").Value = -1089

$ws.Range("H140").Value = 2030.4546
$ws.Range("I140").Value = 2030.4546
$ws.Range("K140").Value = 6091.3638
$ws.Range("M140").Value = -911.3638000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 46000
$ws.Range("J103").Value = 46000
$ws.Range("L103").Value = 46000
$ws.Range("N103").Value = -48344

$ws.Range("H122").Value = 2839.8
$ws.Range("J122").Value = 2844
$ws.Range("L122").Value = 8532
$ws.Range("N122").Value = -13432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 150
$ws.Range("I18").Value = 140
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 140
$ws.Range("L18").Value = 200
$ws.Range("M18").Value = 32
$ws.Range("N18").Value = -544

$ws.Range("H102").Value = 41561.332
$ws.Range("I102").Value = 44839.5
$ws.Range("J102").Value = 39922.25
$ws.Range("K102").Value = 44839.5
$ws.Range("L102").Value = 39922.25
$ws.Range("M102").Value = -41594.5
$ws.Range("N102").Value = -46412.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 31422.25
$ws.Range("J55").Value = 285
$ws.Range("L55").Value = 285
$ws.Range("N55").Value = -839

$ws.Range("H96").Value = 2184
$ws.Range("I96").Value = 1662.909
$ws.Range("J96").Value = 5050
$ws.Range("K96").Value = 1662.909
$ws.Range("L96").Value = 5050
$ws.Range("M96").Value = -289.9090000000001
$ws.Range("N96").Value = -7796

$ws.Range("H132").Value = 9062.741
$ws.Range("I132").Value = 1625.8422
$ws.Range("J132").Value = 26725.375
$ws.Range("K132").Value = 4877.5266
$ws.Range("L132").Value = 80176.125
$ws.Range("M132").Value = -2347.5266
$ws.Range("N132").Value = -85236.125

$ws.Range("H136").Value = 26727.182
$ws.Range("I136").Value = 1999.5
$ws.Range("K136").Value = 5998.5
$ws.Range("M136").Value = -3448.5
